# Insert a new data row at row 88 (pushing existing rows 88:218 down to 89:219)
# and populate it with a new weekly price observation for Coliflor.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(88).Insert()

$ws.Range("A88").Value = 7
$ws.Range("B88").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C88").Value = "Ñuble"
$ws.Range("D88").Value = 44580
$ws.Range("E88").Value = 16
$ws.Range("F88").Value = 100112008
$ws.Range("G88").Value = "Coliflor"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 160
$ws.Range("K88").Value = 850
$ws.Range("L88").Value = 900
$ws.Range("M88").Value = 875
$ws.Range("N88").Value = "$/unidad"
$ws.Range("O88").Value = "Provincia de Diguillín"
$ws.Range("P88").Value = 875
$ws.Range("Q88").Value = 1
$ws.Range("R88").Value = "Hortaliza"
